$d = $word.ActiveDocument

# Locate the placeholder paragraph ("***") that immediately follows the
# "Wi-Fi Networking" heading and precedes the "Security Measures" heading.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq "***") {
        $prev = $d.Paragraphs($i - 1).Range.Text.TrimEnd([char]13)
        if ($prev -eq "Wi-Fi Networking") {
            $target = $i
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate the Wi-Fi Networking placeholder paragraph"
}

# 1) Replace the placeholder text with the first new paragraph's content.
$p = $d.Paragraphs($target)
$p.Range.Text = "It’s important to have wifi coverage throughout the entirety of Greenfields office building, because there are 69 devices that will use wifi and all of them need to stay connected. My recommendation is that the network support should cover for 50% more than the current number of devices, which will be 104 devices. This is to cover for future expansions of the business and prevent issues of device support in the future."

# 2) Add the second new paragraph.
$p = $d.Paragraphs($target)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($target + 1)
$p.Range.Text = "My strategy for planning the wifi coverage is...."

# 3) Add the third new paragraph - two runs: "**DIAGRAM HERE** " and
#    "+ significances of the different colors/numbers***".
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($target + 2)
$p.Range.Text = "**DIAGRAM HERE** "
$insertPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$insertPoint.InsertAfter("+ significances of the different colors/numbers***")
# Force a run boundary between the two sentences (identical formatting would
# otherwise coalesce them into a single run on save).
$p = $d.Paragraphs($target + 2)
$secondRunLen = "+ significances of the different colors/numbers***".Length
$secondRun = $d.Range($p.Range.End - 1 - $secondRunLen, $p.Range.End - 1)
$secondRun.Bold = 1
$secondRun.Bold = 0

# 4) Add the fourth new paragraph.
$p = $d.Paragraphs($target + 2)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($target + 3)
$p.Range.Text = "The network infrastrastructure components that will be needed are..."

# 5) Add the fifth new paragraph.
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($target + 4)
$p.Range.Text = "I recommend using a wireless LAN controller because it offers central management over all the wifi access points (AP) in the building. It also allows updates and changes in configuration to be applied to all access points at once, reducing the possibility that an access point is configured incorrectly. As the IT staff is small, this controller makes the most sense for the business as it streamlines control over the wifi network."

# 6) Add the sixth new paragraph.
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($target + 5)
$p.Range.Text = "The newest and most secure encryption standard for wifi networks is WPA3 (wifi protected access 3). It offers enhanced security features over previous standards and uses SAE (simultaneous authentication of equals) to create keys for connections, making it difficult against dictionary hacking attempts. Another feature is individual data encryption, letting each device have a unique encryption key. So even if the key on 1 device is hacked or compromised, every other devices key on the network will still be protected. As WPA3 is the current IT industry standard of protection, it’s the best encryption method to handle current and new threats to networks."

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
